$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already carry a Synopsis (col H) value - leave them untouched.
# Every other data row (2-32) gets "NA" in column H, using the same
# "blank note" style (s=2) that the old I/J/K filler cells used.
$rowsWithH = @(13, 17, 21, 22, 28, 30, 32)

$ws.Range("I13").Copy() | Out-Null

for ($r = 2; $r -le 32; $r++) {
    if ($rowsWithH -contains $r) {
        continue
    }
    $cell = $ws.Cells.Item($r, 8)
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value = "NA"
}

$excel.CutCopyMode = $false

# The old blank formatting cells in columns I:K are no longer used - drop them
# entirely (this also shrinks the sheet dimension/col widths/row spans down
# to column H).
$ws.Range("I1:K32").Delete() | Out-Null

# Scroll the view over to show column F first and select the now-empty
# columns I:M (mirrors the saved view state in the workbook).
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("I1:M1048576").Select() | Out-Null
